$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work Plan")

function Set-DateCell($range, $year, $month, $day) {
    $d = Get-Date -Year $year -Month $month -Day $day -Hour 0 -Minute 0 -Second 0
    $range.Value = $d
    if ($range.NumberFormat -eq "General") {
        $range.NumberFormat = "m/d/yyyy"
    }
}

# Row 20 - Javascript Calculator (Day 2 actual completion)
Set-DateCell $ws.Range("F20") 2016 7 15
$ws.Range("G20").Value = 100
$ws.Range("H20").Value = "Completed"
$ws.Range("I20").Value = "Second actual commit will be done on the presentation day"
$ws.Rows.Item(20).RowHeight = 30.75

# Row 22 - Developer Check-in code to Repository
Set-DateCell $ws.Range("C22") 2016 7 15
Set-DateCell $ws.Range("D22") 2016 7 15
Set-DateCell $ws.Range("E22") 2016 7 15
Set-DateCell $ws.Range("F22") 2016 7 15
$ws.Range("G22").Value = 100
$ws.Range("H22").Value = "Completed"

# Row 23 - Configure GitLab and Jenkins
Set-DateCell $ws.Range("C23") 2016 7 15
Set-DateCell $ws.Range("D23") 2016 7 15
Set-DateCell $ws.Range("E23") 2016 7 15
Set-DateCell $ws.Range("F23") 2016 7 15
$ws.Range("G23").Value = 100
$ws.Range("H23").Value = "Completed"

# Row 24 - Configure Jenkins and Maven
Set-DateCell $ws.Range("C24") 2016 7 15
Set-DateCell $ws.Range("D24") 2016 7 15
Set-DateCell $ws.Range("E24") 2016 7 15
Set-DateCell $ws.Range("F24") 2016 7 15
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = "Completed"

# Row 25 - Configure Jenkins and Sonarqube (Activity 3 header)
Set-DateCell $ws.Range("C25") 2016 7 18
Set-DateCell $ws.Range("D25") 2016 7 18
Set-DateCell $ws.Range("E25") 2016 7 15
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = "In Progress"

# Row 26 - Deploy App from Ansible to Tomcat
Set-DateCell $ws.Range("C26") 2016 7 18
Set-DateCell $ws.Range("D26") 2016 7 18

# Row 27 - Test app using Selenium
Set-DateCell $ws.Range("C27") 2016 7 18
Set-DateCell $ws.Range("D27") 2016 7 18

# Row 29 - Debug Project Pipeline
Set-DateCell $ws.Range("C29") 2016 7 19
Set-DateCell $ws.Range("D29") 2016 7 19

# Row 30 - Finalize Documents
Set-DateCell $ws.Range("C30") 2016 7 19
Set-DateCell $ws.Range("D30") 2016 7 19

# Update selection to match the final cursor position recorded in the saved file
$ws.Activate()
$ws.Range("E30").Select()
